$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.775216579437256
$ws.Range("B1").Value = 2.185456037521362
$ws.Range("C1").Value = 2.678494453430176
$ws.Range("D1").Value = 6.041849136352539
$ws.Range("E1").Value = 0.854672908782959
